# Weekly update: a new reporting date (2023-03-21, serial 45006) is inserted
# at the top of the data block for this market/product. The two new rows
# (Primera / Segunda quality rows) reuse the same values that used to belong
# to the then-most-recent date, which pushes all the existing rows down by
# two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows right above the (soon to be shifted) block.
$ws.Rows("271:272").Insert()

# The rows that used to be 271:272 are now 273:274 after the insert above.
# Duplicate them into the freshly inserted 271:272 rows, then overwrite the
# date column with the new reporting date.
$ws.Range("A273:R274").Copy() | Out-Null
$ws.Range("A271").PasteSpecial() | Out-Null

$ws.Range("D271:D272").Value2 = 45006

Write-Output "Inserted new reporting rows 271:272 (date 45006) and shifted existing rows down."
